$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: remove the (now stray) _GoBack bookmark that sits after
# "I am planning to add more modules into the program over the time."
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
    Write-Host "Removed old _GoBack bookmark."
}

# ---------------------------------------------------------------------------
# Edit 2: make the "Example: " paragraph bold (paragraph mark + run).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Example:")) {
        $p.Range.Bold = 1
        Write-Host "Bolded paragraph" $i "(Example:)"
    }
}

# ---------------------------------------------------------------------------
# Edit 3: in the LLDP sentence, wrap the first "22" in curly quotes and split
# it out into its own runs, with a fresh _GoBack bookmark placed right after
# the opening quote (i.e. between the quote run and the "22" run).
# ---------------------------------------------------------------------------

# Step 1 - plain text substitution: add the curly quotes around the first
# "22" (keeps everything as one run for now).
$rReplace = $d.Content
[void]$rReplace.Find.Execute("a file name 22. As you can see", $true, $false, $false, $false, $false, $true, 1, $false, "a file name ‘22’. As you can see", 2)

# Step 2 - split the run right before the opening quote, i.e. between
# "...file name " and "‘22’. As you can see...".
$rSplit1 = $d.Content
[void]$rSplit1.Find.Execute("a LLDP check against a file name ")
$rSplit1.Collapse(0)
$d.Bookmarks.Add("TempSplit1", $rSplit1)
$d.Bookmarks.Item("TempSplit1").Delete()

# Step 3 - split right after the opening quote, i.e. between "‘" and "22...".
$rSplit2 = $d.Content
[void]$rSplit2.Find.Execute("against a file name ‘")
$rSplit2.Collapse(0)
$d.Bookmarks.Add("TempSplit2", $rSplit2)
$d.Bookmarks.Item("TempSplit2").Delete()

# Step 4 - split right after "22", i.e. between "22" and "’...".
$rSplit3 = $d.Content
[void]$rSplit3.Find.Execute("name ‘22")
$rSplit3.Collapse(0)
$d.Bookmarks.Add("TempSplit3", $rSplit3)
$d.Bookmarks.Item("TempSplit3").Delete()

# Step 5 - split right after the closing quote, i.e. between "’" and
# ". As you can see...".
$rSplit4 = $d.Content
[void]$rSplit4.Find.Execute("name ‘22’")
$rSplit4.Collapse(0)
$d.Bookmarks.Add("TempSplit4", $rSplit4)
$d.Bookmarks.Item("TempSplit4").Delete()

# Step 6 - drop the new _GoBack bookmark between the opening quote run and
# the "22" run (this is also where the temp split in step 3 already cut).
$rGoBack = $d.Content
[void]$rGoBack.Find.Execute("against a file name ‘")
$rGoBack.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rGoBack)
Write-Host "Inserted new _GoBack bookmark before '22'."

Write-Host "Done."
